$d = $word.ActiveDocument

# Build a minimal OOXML "package" fragment wrapping a single run so it can
# be fed to Range.InsertXML. Using InsertXML (rather than Range.Text /
# Find-Replace) lets us swap just the text-bearing run's contents without
# the engine coalescing/removing the sibling empty <w:r/> runs that sit
# next to several of the paragraphs we need to touch.
function Get-RunPackageXml($runInnerXml) {
    return "<?xml version='1.0'?>" +
        "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>" +
        "<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>" +
        "<pkg:xmlData>" +
        "<w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
        "<w:body><w:p>$runInnerXml</w:p></w:body>" +
        "</w:document>" +
        "</pkg:xmlData></pkg:part></pkg:package>"
}

# Finds the first remaining occurrence of $find in the document and swaps
# the run containing it for a brand-new run holding $replace (optionally
# carrying run formatting given as raw <w:rPr> inner XML in $rpr).
function Replace-RunText($find, $replace, $rpr) {
    $f = $d.Content
    $found = $f.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $find"
        return
    }
    $target = $d.Range($f.Start, $f.End)
    $rprXml = ""
    if ($rpr) { $rprXml = "<w:rPr>$rpr</w:rPr>" }
    $runXml = "<w:r>$rprXml<w:t>$replace</w:t></w:r>"
    $target.InsertXML((Get-RunPackageXml $runXml))
}

# Title (appears twice: Heading1 at top, and a bold run near the end)
Replace-RunText "Play Lines of Magic for Free - A Mysterious Lab Slot Game" "Play Lines of Magic for Free" $null
Replace-RunText "Play Lines of Magic for Free - A Mysterious Lab Slot Game" "Play Lines of Magic for Free" "<w:b/>"

# "What we like" bullet points
Replace-RunText "Simple but engaging gameplay" "Engaging gameplay experience" $null
Replace-RunText "Maximum win amount of 5,000 times bet" "High maximum win amount" $null
Replace-RunText "Mysterious soundtrack perfectly fitting the theme" "Mystical theme and soundtrack" $null

# "What we don't like" bullet points
Replace-RunText "Lacks innovative features" "Lacks animations" $null
Replace-RunText "No elaborate animations" "Not innovative" $null

# Closing italic summary paragraph
Replace-RunText "Experience simple yet engaging gameplay in Lines of Magic, a visually pleasing slot game set in a mysterious laboratory. Play for free and win up to 5,000 times your bet." "Read our review of Lines of Magic and play this engaging slot game for free." "<w:i/>"
